$p = $ppt.ActivePresentation

# --- 1. Update the cached "datetimeFigureOut" field text (6/1/2020 -> 6/5/2020) ---
# This text is cached inside the Date placeholder on the Slide Master and every
# Slide Layout, plus the Notes Master.

# Slide Master's own Date placeholder shape
for ($i = 1; $i -le $p.SlideMaster.Shapes.Count; $i++) {
    $sh = $p.SlideMaster.Shapes.Item($i)
    if ($sh.Name -like "Date*") {
        $sh.TextFrame.TextRange.Text = "6/5/2020"
    }
}

# Every Custom Layout's Date placeholder shape
for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    $lyt = $p.SlideMaster.CustomLayouts.Item($j)
    for ($i = 1; $i -le $lyt.Shapes.Count; $i++) {
        $sh = $lyt.Shapes.Item($i)
        if ($sh.Name -like "Date*") {
            $sh.TextFrame.TextRange.Text = "6/5/2020"
        }
    }
}

# Notes Master date/time header-footer value
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "6/5/2020"

# --- 2. Slide 4 ("Model"): fix up the second bullet under the Content placeholder ---
$slide4 = $p.Slides.Item(4)
for ($i = 1; $i -le $slide4.Shapes.Count; $i++) {
    $sh = $slide4.Shapes.Item($i)
    if ($sh.Name -eq "Content Placeholder 2") {
        $tr = $sh.TextFrame.TextRange
        $para = $tr.Paragraphs(2)
        $para.Runs(1).Text = "Ultimately, a profits a success if we can turn a profit"
    }
}

# --- 3. Slide 5 ("Takeaways"): add a new sub-bullet after the first bullet ---
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $sh = $slide5.Shapes.Item($i)
    if ($sh.Name -eq "Content Placeholder 2") {
        $tr = $sh.TextFrame.TextRange
        $firstPara = $tr.Paragraphs(1)
        $firstPara.InsertAfter("`rAnd we turned a profit!") | Out-Null
        $newPara = $tr.Paragraphs(2)
        $newPara.IndentLevel = 2
    }
}
